$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.115.43'
$ws.Range('E2').Value = '  -0.61%  '

$ws.Range('D3').Value = '1.894.01'
$ws.Range('E3').Value = '  -0.72%  '

$ws.Range('E4').Value = '  +0.19%  '

$ws.Range('D5').Value = '''306.90'
$ws.Range('E5').Value = '  -0.26%  '

$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.10%  '

$ws.Range('D7').Value = '''0.5208'
$ws.Range('E7').Value = '  -0.78%  '

$ws.Range('D8').Value = '''0.3763'
$ws.Range('E8').Value = '  -0.49%  '

$ws.Range('D9').Value = '''0.07275'
$ws.Range('E9').Value = '  +0.14%  '

$ws.Range('E10').Value = '  -0.78%  '

$ws.Range('D11').Value = '''0.9003'
$ws.Range('E11').Value = '  -0.15%  '

$ws.Range('D12').Value = '''0.08218'
$ws.Range('E12').Value = '  -0.52%  '

$ws.Range('D13').Value = '1.962.54'
$ws.Range('E13').Value = '  +2.85%  '

$ws.Range('D14').Value = '''96.42'
$ws.Range('E14').Value = '  +1.23%  '

$ws.Range('D15').Value = '''5.311'
$ws.Range('E15').Value = '  +0.37%  '

$ws.Range('E16').Value = '  +0.22%  '

$ws.Range('D17').Value = '''0.000008614'
$ws.Range('E17').Value = '  +0.07%  '

$ws.Range('D18').Value = '''14.59'
$ws.Range('E18').Value = '  +0.79%  '

$ws.Range('E19').Value = '  +0.28%  '

$ws.Range('D20').Value = '27.143.13'
$ws.Range('E20').Value = '  -0.63%  '

$ws.Range('D21').Value = '''5.081'
$ws.Range('E21').Value = '  +0.20%  '

$ws.Range('E22').Value = '  +0.54%  '

$ws.Range('D23').Value = '''6.420'
$ws.Range('E23').Value = '  -0.65%  '

$ws.Range('B24').Value = 'LidoDAOToken'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D24').Value = '''2.316'
$ws.Range('E24').Value = '  +0.45%  '

$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''148.43'
$ws.Range('E25').Value = '  +1.57%  '

$ws.Range('D26').Value = '''18.20'

$ws.Range('D28').Value = '''115.27'
$ws.Range('E28').Value = '  +0.21%  '

$ws.Range('D29').Value = '''4.809'
$ws.Range('E29').Value = '  -0.12%  '

$ws.Range('D30').Value = '''4.859'
$ws.Range('E30').Value = '  -2.78%  '

$ws.Range('D31').Value = '''0.09206'
$ws.Range('E31').Value = '  -0.24%  '

$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '''0.7940'
$ws.Range('E32').Value = '  -1.56%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '''0.05023'
$ws.Range('E33').Value = '  -0.73%  '

$ws.Range('D34').Value = '''1.217'
$ws.Range('E34').Value = '  -2.08%  '

$ws.Range('E35').Value = '  +2.23%  '

$ws.Range('D36').Value = '''2.957'
$ws.Range('E36').Value = '  -0.36%  '

$ws.Range('D37').Value = '''2.613'
$ws.Range('E37').Value = '  +1.66%  '

$ws.Range('D38').Value = '''0.5726'
$ws.Range('E38').Value = '  -0.07%  '

$ws.Range('D39').Value = '''0.02001'
$ws.Range('E39').Value = '  +0.74%  '

$ws.Range('D40').Value = '''1.076'
$ws.Range('E40').Value = '  +0.01%  '

$ws.Range('D41').Value = '''9.018'
$ws.Range('E41').Value = '  -0.30%  '

$ws.Range('D42').Value = '''6.563'
$ws.Range('E42').Value = '  -0.96%  '

$ws.Range('E43').Value = '  -2.60%  '

$ws.Range('E44').Value = '  -0.17%  '

$ws.Range('D45').Value = '''0.4879'
$ws.Range('E45').Value = '  +0.85%  '

$ws.Range('D46').Value = '''1.001'
$ws.Range('E46').Value = '  +0.06%  '

$ws.Range('E47').Value = '  -1.05%  '

$ws.Range('D48').Value = '''1.621'
$ws.Range('E48').Value = '  +0.16%  '

$ws.Range('D49').Value = '''38.32'
$ws.Range('E49').Value = '  +1.83%  '

$ws.Range('D50').Value = '''63.75'
$ws.Range('E50').Value = '  -0.13%  '

$ws.Range('E51').Value = '  -0.48%  '
